$d = $word.ActiveDocument

# Common run properties (Times New Roman, sz 23) used throughout these edits.
$rPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr>'
$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>'
$pkgFooter = '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- Edit 1 -------------------------------------------------------------
# ". Bundles the MAC address together for delivery to destination node."
# -> three runs: ". Bundles the " / "data with " / "MAC address together for
#    delivery to destination node."
$r1 = $d.Content
$found1 = $r1.Find.Execute(". Bundles the MAC address together for delivery to destination node.")
if ($found1) {
    $target1 = $d.Range($r1.Start, $r1.End)
    $xml1 = $pkgHeader + `
        '<w:r>' + $rPr + '<w:t xml:space="preserve">. Bundles the </w:t></w:r>' + `
        '<w:r>' + $rPr + '<w:t xml:space="preserve">data with </w:t></w:r>' + `
        '<w:r>' + $rPr + '<w:t>MAC address together for delivery to destination node.</w:t></w:r>' + `
        $pkgFooter
    $target1.InsertXML($xml1)
}

# --- Edit 2 -------------------------------------------------------------
# "Multicast: ... first byte must be a 1 for the address to be a multicast.
#  If ... hexadecimal." -> five runs, keeping the trailing
# " Example: 01-00-5E-F4-50-7A" run (and the paragraph's leading <w:tab/>)
# intact/unchanged.
$r2 = $d.Content
$searchText2 = "Multicast: With multicast, a source can send to a group of devices. The low-order bit of the first byte must be a 1 for the address to be a multicast. If the multicast is also an IP multicast, the first 24 bits of the MAC address are 01-00-5E in hexadecimal. Example: 01-00-5E-F4-50-7A"
$found2 = $r2.Find.Execute($searchText2)
if ($found2) {
    # Extend the start back by one character to capture the paragraph's
    # leading tab in the same InsertXML call (avoids a stray/duplicated tab).
    $target2 = $d.Range($r2.Start - 1, $r2.End)
    $xml2 = $pkgHeader + `
        '<w:r>' + $rPr + '<w:tab/><w:t>Multicast: With multicast, a source can send to a group of devices. The low-order bit of the first byte</w:t></w:r>' + `
        '<w:r>' + $rPr + '<w:t xml:space="preserve"> (3 bits)</w:t></w:r>' + `
        '<w:r>' + $rPr + '<w:t xml:space="preserve"> must be a 1 for the address to be a multicast</w:t></w:r>' + `
        '<w:r>' + $rPr + '<w:t xml:space="preserve"> (or 1110)</w:t></w:r>' + `
        '<w:r>' + $rPr + '<w:t>. If the multicast is also an IP multicast, the first 24 bits of the MAC address are 01-00-5E in hexadecimal.</w:t></w:r>' + `
        '<w:r w:rsidR="00B278BA">' + $rPr + '<w:t xml:space="preserve"> Example: 01-00-5E-F4-50-7A</w:t></w:r>' + `
        $pkgFooter
    $target2.InsertXML($xml2)
}

Write-Host "Edit1 found:" $found1 " Edit2 found:" $found2
